$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (updated date 10-11 -> 10-12)
$ws.Name = "Through 2021-10-12"

# Update the October header label
$ws.Range("A11").Value = "October (through 10-12)"

# Row 11 (October) updated values
$ws.Range("B11").Value = 12
$ws.Range("D11").Value = 19
$ws.Range("E11").Value = 30
$ws.Range("F11").Value = 14
$ws.Range("G11").Value = 55
$ws.Range("H11").Value = 79

# Row 12 (Total) updated values
$ws.Range("B12").Value = 238
$ws.Range("D12").Value = 646
$ws.Range("E12").Value = 578
$ws.Range("F12").Value = 436
$ws.Range("G12").Value = 956
$ws.Range("H12").Value = 1329
